# ----------------------------------------------------------------------
# Edit script for algeria_ligue-1_2023-2024.xlsx
# Applies row reordering (matches were re-sorted upstream) and appends
# 4 newly scraped matches (rows 45-48).
# ----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 8 <-> row 9 (columns F:V only; A:E identify the row slot and stay put)
$ws.Range("F8").Value = "Constantine"
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = "MC Alger"
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1.98
$ws.Range("K8").Value = "21/09/2023 05:12"
$ws.Range("L8").Value = 3.68
$ws.Range("M8").Value = "22/09/2023 16:51"
$ws.Range("N8").Value = 2.89
$ws.Range("O8").Value = "21/09/2023 05:12"
$ws.Range("P8").Value = 2.92
$ws.Range("Q8").Value = "22/09/2023 16:51"
$ws.Range("R8").Value = 3.97
$ws.Range("S8").Value = "21/09/2023 05:12"
$ws.Range("T8").Value = 2.26
$ws.Range("U8").Value = "22/09/2023 16:51"
$ws.Range("V8").Value = "https://www.betexplorer.com/football/algeria/ligue-1/constantine-mc-alger/Eoq3MszL/"
$ws.Range("F9").Value = "Ben Aknoun"
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = "ASO Chlef"
$ws.Range("I9").Value = 2
$ws.Range("J9").Value = 2.81
$ws.Range("K9").Value = "22/09/2023 08:13"
$ws.Range("L9").Value = 2.51
$ws.Range("M9").Value = "22/09/2023 16:00"
$ws.Range("N9").Value = 3.04
$ws.Range("O9").Value = "22/09/2023 08:13"
$ws.Range("P9").Value = 3.06
$ws.Range("Q9").Value = "22/09/2023 16:34"
$ws.Range("R9").Value = 2.64
$ws.Range("S9").Value = "22/09/2023 08:13"
$ws.Range("T9").Value = 3.01
$ws.Range("U9").Value = "22/09/2023 16:00"
$ws.Range("V9").Value = "https://www.betexplorer.com/football/algeria/ligue-1/es-ben-aknoun-aso-chlef/WMgbNNKE/"

# Swap row 12 <-> row 13 (columns F:V only; A:E identify the row slot and stay put)
$ws.Range("F12").Value = "Biskra"
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = "US Souf"
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 1.71
$ws.Range("K12").Value = "22/09/2023 15:13"
$ws.Range("L12").Value = 1.65
$ws.Range("M12").Value = "23/09/2023 19:13"
$ws.Range("N12").Value = 3.22
$ws.Range("O12").Value = "22/09/2023 15:13"
$ws.Range("P12").Value = 3.47
$ws.Range("Q12").Value = "23/09/2023 19:13"
$ws.Range("R12").Value = 4.87
$ws.Range("S12").Value = "22/09/2023 15:13"
$ws.Range("T12").Value = 6.17
$ws.Range("U12").Value = "23/09/2023 19:13"
$ws.Range("V12").Value = "https://www.betexplorer.com/football/algeria/ligue-1/biskra-us-souf/KYnDaKS7/"
$ws.Range("F13").Value = "Oran"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = "Saoura"
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 2.54
$ws.Range("K13").Value = "22/09/2023 08:13"
$ws.Range("L13").Value = 2.34
$ws.Range("M13").Value = "23/09/2023 18:03"
$ws.Range("N13").Value = 2.65
$ws.Range("O13").Value = "22/09/2023 08:13"
$ws.Range("P13").Value = 2.71
$ws.Range("Q13").Value = "23/09/2023 18:03"
$ws.Range("R13").Value = 3.11
$ws.Range("S13").Value = "22/09/2023 08:13"
$ws.Range("T13").Value = 3.84
$ws.Range("U13").Value = "23/09/2023 18:03"
$ws.Range("V13").Value = "https://www.betexplorer.com/football/algeria/ligue-1/oran-saoura/YyyeO358/"

# Swap row 23 <-> row 24 (columns F:V only; A:E identify the row slot and stay put)
$ws.Range("F23").Value = "Biskra"
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = "Paradou"
$ws.Range("I23").Value = 5
$ws.Range("J23").Value = 2.18
$ws.Range("K23").Value = "05/10/2023 07:24"
$ws.Range("L23").Value = 2.03
$ws.Range("M23").Value = "06/10/2023 19:34"
$ws.Range("N23").Value = 2.84
$ws.Range("O23").Value = "05/10/2023 07:24"
$ws.Range("P23").Value = 3
$ws.Range("Q23").Value = "06/10/2023 18:05"
$ws.Range("R23").Value = 3.56
$ws.Range("S23").Value = "05/10/2023 07:24"
$ws.Range("T23").Value = 4.33
$ws.Range("U23").Value = "06/10/2023 19:34"
$ws.Range("V23").Value = "https://www.betexplorer.com/football/algeria/ligue-1/biskra-paradou/hhWUzskE/"
$ws.Range("F24").Value = "Oran"
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = "Magra"
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 1.98
$ws.Range("K24").Value = "05/10/2023 07:24"
$ws.Range("L24").Value = 1.65
$ws.Range("M24").Value = "06/10/2023 17:55"
$ws.Range("N24").Value = 2.89
$ws.Range("O24").Value = "05/10/2023 07:24"
$ws.Range("P24").Value = 3.28
$ws.Range("Q24").Value = "06/10/2023 19:03"
$ws.Range("R24").Value = 3.98
$ws.Range("S24").Value = "05/10/2023 07:24"
$ws.Range("T24").Value = 5.91
$ws.Range("U24").Value = "06/10/2023 17:55"
$ws.Range("V24").Value = "https://www.betexplorer.com/football/algeria/ligue-1/oran-magra/WrVYZ04K/"

# Swap row 35 <-> row 36 (columns F:V only; A:E identify the row slot and stay put)
$ws.Range("F35").Value = "Magra"
$ws.Range("G35").Value = 3
$ws.Range("H35").Value = "Ben Aknoun"
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = 1.61
$ws.Range("K35").Value = "11/11/2023 10:12"
$ws.Range("L35").Value = 1.62
$ws.Range("M35").Value = "11/11/2023 14:48"
$ws.Range("N35").Value = 3.51
$ws.Range("O35").Value = "11/11/2023 10:12"
$ws.Range("P35").Value = 3.56
$ws.Range("Q35").Value = "11/11/2023 14:48"
$ws.Range("R35").Value = 5.97
$ws.Range("S35").Value = "11/11/2023 10:12"
$ws.Range("T35").Value = 6.3
$ws.Range("U35").Value = "11/11/2023 14:48"
$ws.Range("V35").Value = "https://www.betexplorer.com/football/algeria/ligue-1/magra-es-ben-aknoun/lCJE0FP6/"
$ws.Range("F36").Value = "Khenchela"
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = "Biskra"
$ws.Range("I36").Value = 1
$ws.Range("J36").Value = 1.74
$ws.Range("K36").Value = "10/11/2023 03:13"
$ws.Range("L36").Value = 1.29
$ws.Range("M36").Value = "11/11/2023 10:24"
$ws.Range("N36").Value = 3.19
$ws.Range("O36").Value = "10/11/2023 03:13"
$ws.Range("P36").Value = 4.87
$ws.Range("Q36").Value = "11/11/2023 14:54"
$ws.Range("R36").Value = 4.9
$ws.Range("S36").Value = "10/11/2023 03:13"
$ws.Range("T36").Value = 13.77
$ws.Range("U36").Value = "11/11/2023 14:54"
$ws.Range("V36").Value = "https://www.betexplorer.com/football/algeria/ligue-1/khenchela-biskra/GbL62yef/"

# Rotate rows 25,26,27 (new25=old27, new26=old25, new27=old26)
$ws.Range("F25").Value = "Ben Aknoun"
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = "US Souf"
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 1.37
$ws.Range("K25").Value = "05/10/2023 09:12"
$ws.Range("L25").Value = 1.84
$ws.Range("M25").Value = "07/10/2023 11:15"
$ws.Range("N25").Value = 4.23
$ws.Range("O25").Value = "05/10/2023 09:12"
$ws.Range("P25").Value = 3.29
$ws.Range("Q25").Value = "07/10/2023 14:50"
$ws.Range("R25").Value = 6.79
$ws.Range("S25").Value = "05/10/2023 09:12"
$ws.Range("T25").Value = 4.65
$ws.Range("U25").Value = "07/10/2023 11:15"
$ws.Range("V25").Value = "https://www.betexplorer.com/football/algeria/ligue-1/es-ben-aknoun-us-souf/Q5UxZKJQ/"
$ws.Range("F26").Value = "CR Belouizdad"
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = "Khenchela"
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 1.46
$ws.Range("K26").Value = "06/10/2023 04:56"
$ws.Range("L26").Value = 1.41
$ws.Range("M26").Value = "07/10/2023 16:01"
$ws.Range("N26").Value = 4.01
$ws.Range("O26").Value = "06/10/2023 04:56"
$ws.Range("P26").Value = 4.17
$ws.Range("Q26").Value = "07/10/2023 16:42"
$ws.Range("R26").Value = 7.72
$ws.Range("S26").Value = "06/10/2023 04:56"
$ws.Range("T26").Value = 9.550000000000001
$ws.Range("U26").Value = "07/10/2023 16:42"
$ws.Range("V26").Value = "https://www.betexplorer.com/football/algeria/ligue-1/cr-belouizdad-khenchela/Q56UOzdJ/"
$ws.Range("F27").Value = "El Bayadh"
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = "ASO Chlef"
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1.79
$ws.Range("K27").Value = "07/10/2023 10:13"
$ws.Range("L27").Value = 2.01
$ws.Range("M27").Value = "07/10/2023 14:14"
$ws.Range("N27").Value = 3.23
$ws.Range("O27").Value = "07/10/2023 10:13"
$ws.Range("P27").Value = 3.1
$ws.Range("Q27").Value = "07/10/2023 16:16"
$ws.Range("R27").Value = 4.85
$ws.Range("S27").Value = "07/10/2023 10:13"
$ws.Range("T27").Value = 4.24
$ws.Range("U27").Value = "07/10/2023 16:40"
$ws.Range("V27").Value = "https://www.betexplorer.com/football/algeria/ligue-1/el-bayadh-aso-chlef/pjB9Tbsl/"

# Append 4 new match rows (45-48), extending dimension to A1:V48.
# First clone formatting from the last existing row (44) so the new rows
# keep the same per-column styles (s="1" on Indice, s="2" on data_partida).
$ws.Range("A44:V44").Copy()
$ws.Range("A45:V45").PasteSpecial(-4122)
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "algeria"
$ws.Range("C45").Value = "ligue-1"
$ws.Range("D45").Value = "2023-2024"
$ws.Range("E45").Value = 45248.625
$ws.Range("F45").Value = "Ben Aknoun"
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = "Paradou"
$ws.Range("I45").Value = 1
$ws.Range("J45").Value = 2.79
$ws.Range("K45").Value = "17/11/2023 01:11"
$ws.Range("L45").Value = 4.23
$ws.Range("M45").Value = "18/11/2023 14:54"
$ws.Range("N45").Value = 2.72
$ws.Range("O45").Value = "17/11/2023 01:11"
$ws.Range("P45").Value = 3.15
$ws.Range("Q45").Value = "18/11/2023 14:56"
$ws.Range("R45").Value = 2.67
$ws.Range("S45").Value = "17/11/2023 01:11"
$ws.Range("T45").Value = 1.99
$ws.Range("U45").Value = "18/11/2023 14:54"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/algeria/ligue-1/es-ben-aknoun-paradou/QBmdDh9C/"

$ws.Range("A44:V44").Copy()
$ws.Range("A46:V46").PasteSpecial(-4122)
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "algeria"
$ws.Range("C46").Value = "ligue-1"
$ws.Range("D46").Value = "2023-2024"
$ws.Range("E46").Value = 45248.70833333334
$ws.Range("F46").Value = "ASO Chlef"
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = "ES Setif"
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 1.73
$ws.Range("K46").Value = "17/11/2023 01:11"
$ws.Range("L46").Value = 1.78
$ws.Range("M46").Value = "18/11/2023 16:56"
$ws.Range("N46").Value = 3.19
$ws.Range("O46").Value = "17/11/2023 01:11"
$ws.Range("P46").Value = 3
$ws.Range("Q46").Value = "18/11/2023 16:58"
$ws.Range("R46").Value = 4.76
$ws.Range("S46").Value = "17/11/2023 01:11"
$ws.Range("T46").Value = 4.33
$ws.Range("U46").Value = "18/11/2023 16:58"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/algeria/ligue-1/aso-chlef-es-setif/KYn0CCOI/"

$ws.Range("A44:V44").Copy()
$ws.Range("A47:V47").PasteSpecial(-4122)
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "algeria"
$ws.Range("C47").Value = "ligue-1"
$ws.Range("D47").Value = "2023-2024"
$ws.Range("E47").Value = 45248.75
$ws.Range("F47").Value = "Biskra"
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = "USM Alger"
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 2.02
$ws.Range("K47").Value = "17/11/2023 01:11"
$ws.Range("L47").Value = 2.67
$ws.Range("M47").Value = "18/11/2023 17:24"
$ws.Range("N47").Value = 2.93
$ws.Range("O47").Value = "17/11/2023 01:11"
$ws.Range("P47").Value = 2.99
$ws.Range("Q47").Value = "18/11/2023 16:09"
$ws.Range("R47").Value = 3.79
$ws.Range("S47").Value = "17/11/2023 01:11"
$ws.Range("T47").Value = 2.87
$ws.Range("U47").Value = "18/11/2023 17:24"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/algeria/ligue-1/biskra-usm-alger/bPo4BWwP/"

$ws.Range("A44:V44").Copy()
$ws.Range("A48:V48").PasteSpecial(-4122)
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "algeria"
$ws.Range("C48").Value = "ligue-1"
$ws.Range("D48").Value = "2023-2024"
$ws.Range("E48").Value = 45249.79166666666
$ws.Range("F48").Value = "CR Belouizdad"
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = "Kabylie"
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 1.86
$ws.Range("K48").Value = "18/11/2023 07:12"
$ws.Range("L48").Value = 1.61
$ws.Range("M48").Value = "19/11/2023 18:32"
$ws.Range("N48").Value = 3
$ws.Range("O48").Value = "18/11/2023 07:12"
$ws.Range("P48").Value = 3.5
$ws.Range("Q48").Value = "19/11/2023 18:32"
$ws.Range("R48").Value = 4.32
$ws.Range("S48").Value = "18/11/2023 07:12"
$ws.Range("T48").Value = 6.71
$ws.Range("U48").Value = "19/11/2023 18:32"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/algeria/ligue-1/cr-belouizdad-kabylie/xAHpKOsP/"

